$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

$ws.Range("D2").Value = -0.0653
$ws.Range("D3").Value = -0.0653
$ws.Range("G2").Value = 0.1563593932322054
$ws.Range("G3").Value = 0.1563593932322054
$ws.Range("H2").Value = 0.1563593932322054
$ws.Range("H3").Value = 0.1563593932322054
$ws.Range("I2").Value = -0.02023604279761404
$ws.Range("I3").Value = -0.02023604279761404
$ws.Range("J2").Value = -0.02023604279761404
$ws.Range("J3").Value = -0.02023604279761404
$ws.Range("K2").Value = -182.2
$ws.Range("K3").Value = -182.2
$ws.Range("L2").Value = -0.0966373183409356
$ws.Range("L3").Value = -0.0966373183409356
$ws.Range("M2").Value = 0.002
$ws.Range("M3").Value = 0.002
$ws.Range("N2").Value = 0.000002001801621459313
$ws.Range("N3").Value = 0.000002001801621459313
$ws.Range("O2").Value = -0.00001097694840834248
$ws.Range("O3").Value = -0.00001097694840834248
$ws.Range("P2").Value = 0.002
$ws.Range("P3").Value = 0.002
$ws.Range("Q2").Value = 0.000002001801621459313
$ws.Range("Q3").Value = 0.000002001801621459313
$ws.Range("R2").Value = -0.00001097694840834248
$ws.Range("R3").Value = -0.00001097694840834248
$ws.Range("U2").Value = 149.9
$ws.Range("U3").Value = 149.9
$ws.Range("V2").Value = 0.1500350315283756
$ws.Range("V3").Value = 0.1500350315283756
$ws.Range("W2").Value = -0.1258808898714937
$ws.Range("W3").Value = -0.1258808898714937
$ws.Range("X2").Value = 0.06847942397093404
$ws.Range("X3").Value = 0.06847942397093404
$ws.Range("Y2").Value = -0.1943603138424277
$ws.Range("Y3").Value = -0.1943603138424277
$ws.Range("Z2").Value = 0.7996894638093194
$ws.Range("Z3").Value = 0.7996894638093194
$ws.Range("AA2").Value = -0.01618255021444641
$ws.Range("AA3").Value = -0.01618255021444641
$ws.Range("AB2").Value = 0.04543403971331211
$ws.Range("AB3").Value = 0.04543403971331211
$ws.Range("AC2").Value = -0.06161658992775852
$ws.Range("AC3").Value = -0.06161658992775852
$ws.Range("AD2").Value = 1363.1
$ws.Range("AD3").Value = 1363.1
$ws.Range("AE2").Value = 5.565175453107585
$ws.Range("AE3").Value = 5.565175453107585
$ws.Range("AF2").Value = 1368.665175453107
$ws.Range("AF3").Value = 1368.665175453107
$ws.Range("AG2").Value = 1218.765175453107
$ws.Range("AG3").Value = 1218.765175453107
$ws.Range("AH2").Value = 0.578040926373196
$ws.Range("AH3").Value = 0.578040926373196
$ws.Range("AI2").Value = 0.5365914697957068
$ws.Range("AI3").Value = 0.5365914697957068
$ws.Range("AJ2").Value = 0.5495217603586363
$ws.Range("AJ3").Value = 0.5495217603586363
$ws.Range("AK2").Value = 0.5076569703336705
$ws.Range("AK3").Value = 0.5076569703336705
$ws.Range("AL2").Value = 74.59999999999999
$ws.Range("AL3").Value = 74.59999999999999
$ws.Range("AM2").Value = 57.7
$ws.Range("AM3").Value = 57.7
$ws.Range("AN2").Value = 5.327522864066285
$ws.Range("AN3").Value = 5.327522864066285
$ws.Range("AO2").Value = -0.5415549597855228
$ws.Range("AO3").Value = -0.5415549597855228
$ws.Range("AP2").Value = 4.763406454518515
$ws.Range("AP3").Value = 4.763406454518515
$ws.Range("AQ2").Value = -0.7001733102253033
$ws.Range("AQ3").Value = -0.7001733102253033
